$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 used to hold the single "Accuracy over PyType" stat in E29/F29. A new
# "Scalpel Accuracy:" stat is inserted into C29/D29, and the old
# "Accuracy over PyType" stat is pushed down into a brand new row 30 (E30/F30).

# --- Row 29: add the new "Scalpel Accuracy:" label + value, clear old E/F ---
$ws.Range("C29").Value = "Scalpel Accuracy:"
$ws.Range("D29").Value = 766.67
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = ""

# --- Row 30 (new): carry the "Accuracy over PyType" label + value down ---
$ws.Range("E30").Value = "Accuracy over PyType"
$ws.Range("F30").Value = 33.33

# Give every cell in the new row 30 (and keep row 29) the same fill/style used
# throughout the rest of the data rows (style index "2" -- solid white fill),
# by mirroring the existing row 28/29 formatting onto the untouched A30:D30 cells.
$ws.Range("A30:F30").Interior.Color = $ws.Range("A29").Interior.Color
